# This patch removes internal rounding of floating point values used
# to compute the figures in this reference sheet. The previously
# rounded/integer values in B2:E13 are replaced with their unrounded
# floating point equivalents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(810.6731951491925,1157.387157691246,1354.354347540363,1505.730094935563),
    @(853.6162607233118,1208.207581561454,1401.132457741924,1543.621119192856),
    @(734.5948601966046,1063.700906389054,1261.470515273808,1421.612905331771),
    @(915.9270000733128,1281.062412258313,1454.593041708381,1599.557046887749),
    @(900.2329883775661,1262.312376596757,1437.299634013683,1585.162946869089),
    @(938.8517375560593,1305.139994837836,1488.506807391738,1612.1304178267),
    @(833.1799429183484,1187.049016695813,1396.223082153899,1524.199078212293),
    @(945.7389134525196,1317.148662192248,1501.802711532365,1625.195807347984),
    @(961.3154436913049,1315.913788041324,1473.909083247515,1570.125201648339),
    @(952.5553249171585,1302.951897844415,1459.287229046187,1555.234797529095),
    @(874.7493301769194,1157.343514172884,1251.385818327302,1297.715867080148),
    @(956.7539180159872,1308.382254672927,1464.790387914299,1560.211442101217)
)

$cols = @("B", "C", "D", "E")
$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $row = $startRow + $i
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range("$($cols[$j])$row").Value = $rowValues[$j]
    }
}
